$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns to use the format-version suffix instead of old/new ---
$cols = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $cols[$i] + "_FV2210"
}
# Column K (11) stays "diff"
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $cols[$i] + "_FV2304"
}

# --- 2. Stash the existing header formatting so we can restore it after the table is created ---
# (ListObjects.Add() otherwise bakes the current header format into a header-row dxf,
#  which the source workbook does not have)
$helper = $ws.Range("A200:U200")
$helper.Value = "x"
$hdr = $ws.Range("A1:U1")
$hdr.Copy()
$helper.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$hdr.ClearFormats()

# --- 3. Turn the used range into an Excel Table ---
$rng = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 4. Restore the original header formatting ---
$helper.Copy()
$hdr.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$helper.Clear()

# --- 5. Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
